# Add "First Letter of Cat's Name" and "Cat Age (yr)" columns to the Data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Header row
$ws.Range("D1").Value = "First Letter of Cat's Name"
$ws.Range("E1").Value = "Cat Age (yr)"
$ws.Range("D1:E1").Font.Bold = $true

# Data rows: first letter of cat's name, cat age (yr)
$catLetter = @("A","A","Z","D","B","D","P","T","M","G","M","E","C","R")
$catAge = @(5,7,3,2,4,9,2,10,15,3,3,17,4,1)

for ($i = 0; $i -lt $catLetter.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $catLetter[$i]
    $ws.Cells.Item($row, 5).Value = $catAge[$i]
}
